$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.799.99"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "3.798.48"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "701.18"
$ws.Range("E5").Value = "  +5.91%  "
$ws.Range("D6").Value = "172.46"
$ws.Range("E6").Value = "  +3.56%  "
$ws.Range("D7").Value = "3.798.31"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").Value = "7.25"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E13").Value = "  +6.08%  "
$ws.Range("D14").Value = "36.03"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "4.439.11"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "3.794.43"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "70.729.28"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "10.96"
$ws.Range("E21").Value = "  +13.91%  "
$ws.Range("D22").Value = "480.53"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "84.06"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").Value = "0.0000143"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "12.30"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "10.49"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "3.948.35"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "3.13"
$ws.Range("E31").Value = "  +12.52%  "
$ws.Range("E32").Value = "  +4.41%  "
$ws.Range("D33").Value = "2.29"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +8.16%  "
$ws.Range("D35").Value = "29.41"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").Value = "9.25"
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").Value = "2.22"
$ws.Range("E41").Value = "  +9.50%  "
$ws.Range("D42").Value = "0.980"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D45").Value = "0.000316"
$ws.Range("E45").Value = "  +15.84%  "
$ws.Range("D46").Value = "163.49"
$ws.Range("E46").Value = "  +3.58%  "
$ws.Range("D47").Value = "48.85"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").Value = "44.35"
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("D49").Value = "1.40"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "8.61"
$ws.Range("E51").Value = "  +1.83%  "

$ws.Range("D2:D51").Style = "Normal"
